$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.639.26"
$ws.Range("E2").Value = "  +2.72%  "
$ws.Range("D3").Value = "2.003.79"
$ws.Range("E3").Value = "  +6.50%  "
$ws.Range("D4").Value = "'1.007"
$ws.Range("E4").Value = "  +0.28%  "
$ws.Range("D5").Value = "'328.81"
$ws.Range("E5").Value = "  +1.03%  "
$ws.Range("D6").Value = "'1.005"
$ws.Range("E6").Value = "  +0.25%  "
$ws.Range("D7").Value = "'0.4720"
$ws.Range("E7").Value = "  +2.81%  "
$ws.Range("D8").Value = "'0.3965"
$ws.Range("E8").Value = "  +2.10%  "
$ws.Range("D9").Value = "'47.00"
$ws.Range("E9").Value = "  +1.05%  "
$ws.Range("D10").Value = "'0.07967"
$ws.Range("E10").Value = "  +1.45%  "
$ws.Range("D11").Value = "'1.007"
$ws.Range("E11").Value = "  +2.20%  "
$ws.Range("D12").Value = "'22.74"
$ws.Range("E12").Value = "  +4.49%  "
$ws.Range("D13").Value = "2.003.61"
$ws.Range("E13").Value = "  +3.34%  "
$ws.Range("D14").Value = "'7.271"
$ws.Range("E14").Value = "  +3.71%  "
$ws.Range("D15").Value = "'5.895"
$ws.Range("E15").Value = "  +4.06%  "
$ws.Range("D16").Value = "'0.07164"
$ws.Range("E16").Value = "  +3.15%  "
$ws.Range("D17").Value = "'89.21"
$ws.Range("E17").Value = "  +1.10%  "
$ws.Range("D18").Value = "'1.008"
$ws.Range("E18").Value = "  +0.46%  "
$ws.Range("D19").Value = "'0.00001002"
$ws.Range("E19").Value = "  +0.56%  "
$ws.Range("D20").Value = "'17.42"
$ws.Range("E20").Value = "  +2.59%  "
$ws.Range("D21").Value = "'1.006"
$ws.Range("E21").Value = "  +0.41%  "
$ws.Range("D22").Value = "29.738.76"
$ws.Range("E22").Value = "  +2.91%  "
$ws.Range("D23").Value = "'5.556"
$ws.Range("E23").Value = "  +5.41%  "
$ws.Range("D24").Value = "'11.29"
$ws.Range("E24").Value = "  +3.14%  "
$ws.Range("D25").Value = "2.275.92"
$ws.Range("E25").Value = "  +5.38%  "
$ws.Range("D26").Value = "'2.135"
$ws.Range("E26").Value = "  +2.26%  "
$ws.Range("D27").Value = "'158.78"
$ws.Range("E27").Value = "  +2.17%  "
$ws.Range("D28").Value = "'19.76"
$ws.Range("E28").Value = "  +2.48%  "
$ws.Range("D29").Value = "'5.981"
$ws.Range("E29").Value = "  +0.48%  "
$ws.Range("D30").Value = "'120.79"
$ws.Range("E30").Value = "  +2.82%  "
$ws.Range("D31").Value = "'1.970"
$ws.Range("E31").Value = "  +2.15%  "
$ws.Range("D32").Value = "'0.09495"
$ws.Range("E32").Value = "  +1.63%  "
$ws.Range("D33").Value = "'0.8988"
$ws.Range("E33").Value = "  -0.34%  "
$ws.Range("D34").Value = "'5.307"
$ws.Range("E34").Value = "  +0.62%  "
$ws.Range("D35").Value = "'1.345"
$ws.Range("E35").Value = "  +1.39%  "
$ws.Range("D36").Value = "'3.195"
$ws.Range("E36").Value = "  -2.19%  "
$ws.Range("D37").Value = "'0.05850"
$ws.Range("E37").Value = "  +1.58%  "
$ws.Range("B38").Value = "PEPE"
$ws.Range("C38").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D38").Value = "'0.000003434"
$ws.Range("E38").Value = "  +111.11%  "
$ws.Range("B39").Value = "TrustWalletToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D39").Value = "'1.180"
$ws.Range("E39").Value = "  -0.52%  "
$ws.Range("B40").Value = "VeChain"
$ws.Range("C40").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D40").Value = "'0.02140"
$ws.Range("E40").Value = "  +3.34%  "
$ws.Range("D41").Value = "'7.932"
$ws.Range("E41").Value = "  +3.66%  "
$ws.Range("B42").Value = "TheSandbox"
$ws.Range("C42").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D42").Value = "'0.5784"
$ws.Range("E42").Value = "  +2.25%  "
$ws.Range("B43").Value = "Algorand"
$ws.Range("C43").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D43").Value = "'0.1828"
$ws.Range("E43").Value = "  +3.55%  "
$ws.Range("B44").Value = "Aptos"
$ws.Range("C44").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D44").Value = "'9.892"
$ws.Range("E44").Value = "  +2.22%  "
$ws.Range("B45").Value = "EnergySwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D45").Value = "'12.22"
$ws.Range("E45").Value = "  +2.44%  "
$ws.Range("B46").Value = "Decentraland"
$ws.Range("C46").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D46").Value = "'0.5404"
$ws.Range("E46").Value = "  +0.99%  "
$ws.Range("B47").Value = "MXToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D47").Value = "'2.659"
$ws.Range("E47").Value = "  +4.81%  "
$ws.Range("B48").Value = "RenderToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D48").Value = "'2.163"
$ws.Range("E48").Value = "  -4.25%  "
$ws.Range("B49").Value = "Cronos"
$ws.Range("C49").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D49").Value = "'0.06998"
$ws.Range("E49").Value = "  -0.57%  "
$ws.Range("B50").Value = "NEARProtocol"
$ws.Range("C50").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D50").Value = "'1.878"
$ws.Range("E50").Value = "  +1.61%  "
$ws.Range("B51").Value = "Quant"
$ws.Range("C51").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D51").Value = "'114.73"
$ws.Range("E51").Value = "  +1.71%  "
